# Re-shuffle ride assignments across rows 3-20 of the "Ride Assignments" sheet.
# Each block below sets a cell's display name; where the fill color (style)
# also changes, format is first copied from a stable same-style anchor cell
# (B2=pink/style1, B3=green/style6, B5=blue/style7, B4=cream/style9) via
# Copy+PasteSpecial(xlPasteFormats) so the existing style index is reused
# instead of Excel minting a near-duplicate style entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 'Aaron duong'
$ws.Range("B2").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E3").Value = 'Ethan Yu'
$ws.Range("F3").Value = 'Jiwang Lee'
$ws.Range("G3").Value = 'Gabriel Ni'
$ws.Range("B5").Copy()
$ws.Range("L3").PasteSpecial(-4122)
$ws.Range("L3").Value = 'derek liang '
$ws.Range("B3").Copy()
$ws.Range("M3").PasteSpecial(-4122)
$ws.Range("M3").Value = 'Darius Ajebon '
$ws.Range("B2").Copy()
$ws.Range("N3").PasteSpecial(-4122)
$ws.Range("N3").Value = 'Christina Ko'
$ws.Range("O3").Value = 'Hannah Kim'

# Row 4
$ws.Range("D4").Value = 'Ella Lu'
$ws.Range("B2").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("E4").Value = 'Joanna Wei'
$ws.Range("F4").Value = 'Israel Haile'
$ws.Range("G4").Value = 'Lucy Han'
$ws.Range("K4").Value = 'Israel Haile'
$ws.Range("B5").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("L4").Value = 'Daniel Kuo'
$ws.Range("B3").Copy()
$ws.Range("M4").PasteSpecial(-4122)
$ws.Range("M4").Value = 'Lucy Han'
$ws.Range("B2").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N4").Value = 'Sehyun Jung'
$ws.Range("O4").Value = 'Kyle Hwang'

# Row 5
$ws.Range("C5").Value = 'Grace Park'
$ws.Range("D5").Value = 'Maya Habraken '
$ws.Range("B2").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("E5").Value = 'Sehyun Jung'
$ws.Range("G5").Value = 'Daniel Kim '
$ws.Range("K5").Value = 'Rachel Kim'
$ws.Range("B5").Copy()
$ws.Range("L5").PasteSpecial(-4122)
$ws.Range("L5").Value = 'Taeho Choe'
$ws.Range("B3").Copy()
$ws.Range("M5").PasteSpecial(-4122)
$ws.Range("M5").Value = 'Emily Yang'
$ws.Range("B2").Copy()
$ws.Range("N5").PasteSpecial(-4122)
$ws.Range("N5").Value = 'Grace Park'

# Row 6
$ws.Range("C6").Value = 'Nathanael Wang'
$ws.Range("D6").Value = 'Grace Kwon'
$ws.Range("B2").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("E6").Value = 'JJ Lee'
$ws.Range("F6").Value = 'Taeho Choe'
$ws.Range("G6").Value = 'Daniel Song'
$ws.Range("K6").Value = 'Aaron duong'
$ws.Range("B5").Copy()
$ws.Range("L6").PasteSpecial(-4122)
$ws.Range("L6").Value = 'Grace Sowon Park '
$ws.Range("B3").Copy()
$ws.Range("M6").PasteSpecial(-4122)
$ws.Range("M6").Value = 'Daniel Kim '
$ws.Range("B2").Copy()
$ws.Range("N6").PasteSpecial(-4122)
$ws.Range("N6").Value = 'Nathanael Wang'

# Row 10
$ws.Range("B4").Copy()
$ws.Range("K10").PasteSpecial(-4122)
$ws.Range("K10").Value = 'April Tong'
$ws.Range("B5").Copy()
$ws.Range("L10").PasteSpecial(-4122)
$ws.Range("L10").Value = 'Joann Jung'
$ws.Range("O10").Value = 'Grace Kwon'

# Row 11
$ws.Range("E11").Value = 'Faith Chen'
$ws.Range("B4").Copy()
$ws.Range("K11").PasteSpecial(-4122)
$ws.Range("K11").Value = 'Pedro Flores-Teran'
$ws.Range("B5").Copy()
$ws.Range("L11").PasteSpecial(-4122)
$ws.Range("L11").Value = 'Hyeongjun Son'
$ws.Range("O11").Value = 'Benjamin Kim'

# Row 12
$ws.Range("G12").Value = 'Austin Lee'
$ws.Range("B4").Copy()
$ws.Range("K12").PasteSpecial(-4122)
$ws.Range("K12").Value = 'Shayla Nguyen'
$ws.Range("B5").Copy()
$ws.Range("L12").PasteSpecial(-4122)
$ws.Range("L12").Value = 'Ella Lu'
$ws.Range("O12").Value = 'Hannah Zhang'

# Row 13
$ws.Range("G13").Value = '김예림'
$ws.Range("B4").Copy()
$ws.Range("K13").PasteSpecial(-4122)
$ws.Range("K13").Value = 'Sam Ko'
$ws.Range("B5").Copy()
$ws.Range("L13").PasteSpecial(-4122)
$ws.Range("L13").Value = 'Ella'
$ws.Range("O13").Value = 'Maya Habraken '

# Row 17
$ws.Range("C17").Value = 'Ella'
$ws.Range("D17").Value = 'Emily Yang'
$ws.Range("E17").Value = 'Hyeongjun Son'
$ws.Range("B5").Copy()
$ws.Range("F17").PasteSpecial(-4122)
$ws.Range("F17").Value = 'Elie Park'
$ws.Range("K17").Value = 'Ethan Yu'
$ws.Range("B4").Copy()
$ws.Range("P17").PasteSpecial(-4122)
$ws.Range("P17").Value = 'Jane Yoo (Back home 💙)'

# Row 18
$ws.Range("C18").Value = 'Hannah Kim'
$ws.Range("D18").Value = 'Jocelyn Youn'
$ws.Range("E18").Value = 'Joann Jung'
$ws.Range("B5").Copy()
$ws.Range("F18").PasteSpecial(-4122)
$ws.Range("F18").Value = 'Grace Sowon Park '
$ws.Range("K18").Value = 'JJ Lee'
$ws.Range("B3").Copy()
$ws.Range("P18").PasteSpecial(-4122)
$ws.Range("P18").Value = 'Daniel Song (Back home 💙)'

# Row 19
$ws.Range("C19").Value = 'Jeffery Huang'
$ws.Range("D19").Value = 'helena song'
$ws.Range("E19").Value = 'Daniel Kuo'
$ws.Range("B5").Copy()
$ws.Range("F19").PasteSpecial(-4122)
$ws.Range("F19").Value = 'Kyle Hwang'

# Row 20
$ws.Range("C20").Value = 'derek liang '
$ws.Range("E20").Value = 'Hannah Zhang'
$ws.Range("B5").Copy()
$ws.Range("F20").PasteSpecial(-4122)
$ws.Range("F20").Value = 'Benjamin Kim'

# Release the marching-ants clipboard selection left by Copy().
$excel.CutCopyMode = $false
